$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.73 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2: 0.27 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3: 97 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4: 361 -> 1443
$t.Cell(4,1).Range.Text = "1443"

# Row 6: 0.00035 -> 0.00059
$t.Cell(6,1).Range.Text = "0.00059"

# Row 7: 0.00008 -> 0.00019
$t.Cell(7,1).Range.Text = "0.00019"

# Row 8: 0.00002 -> 0.00005
$t.Cell(8,1).Range.Text = "0.00005"

# Row 9: 0.00006 -> 0.00028
$t.Cell(9,1).Range.Text = "0.00028"

# Row 10: 0.00007 -> 0.00033
$t.Cell(10,1).Range.Text = "0.00033"

# Row 11: 0.00009 -> 0.00041
$t.Cell(11,1).Range.Text = "0.00041"

# Row 12: 0.02804 -> 0.26833
$t.Cell(12,1).Range.Text = "0.26833"

# Row 44: collapse multi-run tabbed stats line down to 99.73
$t.Cell(44,1).Range.Text = "99.73"

# Row 45: collapse multi-run tabbed stats line down to 0.27
$t.Cell(45,1).Range.Text = "0.27"

# Row 46: collapse multi-run tabbed stats line down to 97
$t.Cell(46,1).Range.Text = "97"
